$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 3.272327238179451
$ws.Range("C2").Value = 1.626987699542094
$ws.Range("D2").Value = 0.1496068669990043
$ws.Range("E2").Value = 0.5333859586016987
$ws.Range("G2").Value = 5.582307763322248

$ws.Range("B3").Value = 1.445647641019636
$ws.Range("C3").Value = 1.626987699542094
$ws.Range("D3").Value = 0.7210945179870265
$ws.Range("E3").Value = 0.5333859586016987
$ws.Range("G3").Value = 4.327115817150455

$ws.Range("B4").Value = 3.272327238179451
$ws.Range("C4").Value = 1.626987699542094
$ws.Range("D4").Value = 0.7210945179870265
$ws.Range("E4").Value = 0.5333859586016987
$ws.Range("G4").Value = 6.15379541431027

$ws.Range("B5").Value = 3.272327238179451
$ws.Range("C5").Value = 1.626987699542094
$ws.Range("D5").Value = 0.1496068669990043
$ws.Range("E5").Value = 0.5333859586016987
$ws.Range("G5").Value = 5.582307763322248

$ws.Range("B6").Value = 3.272327238179451
$ws.Range("C6").Value = 1.626987699542094
$ws.Range("D6").Value = 0.7210945179870265
$ws.Range("E6").Value = 13.86384647080068
$ws.Range("G6").Value = 19.48425592650926

$ws.Range("B7").Value = 0.1169995834814548
$ws.Range("C7").Value = 0.3048912486333797
$ws.Range("D7").Value = 3.223369029078222
$ws.Range("E7").Value = 0.5333859586016987
$ws.Range("G7").Value = 4.178645819794754

$ws.Range("B8").Value = 1.445647641019636
$ws.Range("C8").Value = 1.626987699542094
$ws.Range("D8").Value = 0.7210945179870265
$ws.Range("E8").Value = 0.5333859586016987
$ws.Range("G8").Value = 4.327115817150455

$ws.Range("B9").Value = 0.0000119090754144846
$ws.Range("C9").Value = 0.002658071450198252
$ws.Range("D9").Value = 0.7210945179870265
$ws.Range("E9").Value = 0.5333859586016987
$ws.Range("G9").Value = 1.257150457114338

$ws.Range("B10").Value = 0.6545652718822623
$ws.Range("C10").Value = 1.626987699542094
$ws.Range("D10").Value = 0.1496068669990043
$ws.Range("E10").Value = 0.5333859586016987
$ws.Range("G10").Value = 2.964545797025059

$ws.Range("B11").Value = 3.272327238179451
$ws.Range("C11").Value = 1.626987699542094
$ws.Range("D11").Value = 3.223369029078222
$ws.Range("E11").Value = 0.5333859586016987
$ws.Range("G11").Value = 8.656069925401464

$ws.Range("B12").Value = 0.04172184405617529
$ws.Range("C12").Value = 0.04103571897497393
$ws.Range("D12").Value = 3.223369029078222
$ws.Range("E12").Value = 0.5333859586016987
$ws.Range("G12").Value = 3.83951255071107

$ws.Range("B13").Value = 1.445647641019636
$ws.Range("C13").Value = 1.626987699542094
$ws.Range("D13").Value = 0.1496068669990043
$ws.Range("E13").Value = 0.5333859586016987
$ws.Range("G13").Value = 3.755628166162433

$ws.Range("B14").Value = 3.272327238179451
$ws.Range("C14").Value = 0.3048912486333797
$ws.Range("D14").Value = 18.71679738969934
$ws.Range("E14").Value = 13.86384647080068
$ws.Range("G14").Value = 36.15786234731286

$ws.Range("B15").Value = 1.445647641019636
$ws.Range("C15").Value = 1.626987699542094
$ws.Range("D15").Value = 0.7210945179870265
$ws.Range("E15").Value = 0.5333859586016987
$ws.Range("G15").Value = 4.327115817150455

$ws.Range("B16").Value = 3.272327238179451
$ws.Range("C16").Value = 1.626987699542094
$ws.Range("D16").Value = 0.1496068669990043
$ws.Range("E16").Value = 0.5333859586016987
$ws.Range("G16").Value = 5.582307763322248

$ws.Range("B17").Value = 1.445647641019636
$ws.Range("C17").Value = 1.626987699542094
$ws.Range("D17").Value = 3.223369029078222
$ws.Range("E17").Value = 0.5333859586016987
$ws.Range("G17").Value = 6.82939032824165

$ws.Range("B18").Value = 0.1169995834814548
$ws.Range("C18").Value = 0.3048912486333797
$ws.Range("D18").Value = 0.7210945179870265
$ws.Range("E18").Value = 0.5333859586016987
$ws.Range("G18").Value = 1.67637130870356

$ws.Range("B19").Value = 3.272327238179451
$ws.Range("C19").Value = 1.626987699542094
$ws.Range("D19").Value = 0.1496068669990043
$ws.Range("E19").Value = 0.5333859586016987
$ws.Range("G19").Value = 5.582307763322248

$ws.Range("B20").Value = 0.1169995834814548
$ws.Range("C20").Value = 0.04103571897497393
$ws.Range("D20").Value = 0.7210945179870265
$ws.Range("E20").Value = 0.5333859586016987
$ws.Range("G20").Value = 1.412515779045154

$ws.Range("B21").Value = 0.04172184405617529
$ws.Range("C21").Value = 0.04103571897497393
$ws.Range("D21").Value = 18.71679738969934
$ws.Range("E21").Value = 0.5333859586016987
$ws.Range("G21").Value = 19.33294091133218

$ws.Range("B22").Value = 3.272327238179451
$ws.Range("C22").Value = 1.626987699542094
$ws.Range("D22").Value = 3.223369029078222
$ws.Range("E22").Value = 0.5333859586016987
$ws.Range("G22").Value = 8.656069925401464

$ws.Range("B23").Value = 3.272327238179451
$ws.Range("C23").Value = 1.626987699542094
$ws.Range("D23").Value = 0.1496068669990043
$ws.Range("E23").Value = 0.5333859586016987
$ws.Range("G23").Value = 5.582307763322248

$ws.Range("B24").Value = 3.272327238179451
$ws.Range("C24").Value = 1.626987699542094
$ws.Range("D24").Value = 0.7210945179870265
$ws.Range("E24").Value = 0.5333859586016987
$ws.Range("G24").Value = 6.15379541431027

$ws.Range("B25").Value = 0.00009552326474482342
$ws.Range("C25").Value = 0.04103571897497393
$ws.Range("D25").Value = 0.7210945179870265
$ws.Range("E25").Value = 0.5333859586016987
$ws.Range("G25").Value = 1.295611718828444

$ws.Range("B26").Value = 3.272327238179451
$ws.Range("C26").Value = 1.626987699542094
$ws.Range("D26").Value = 3.223369029078222
$ws.Range("E26").Value = 0.5333859586016987
$ws.Range("G26").Value = 8.656069925401464

$ws.Range("B27").Value = 1.445647641019636
$ws.Range("C27").Value = 1.626987699542094
$ws.Range("D27").Value = 0.7210945179870265
$ws.Range("E27").Value = 0.5333859586016987
$ws.Range("G27").Value = 4.327115817150455

$ws.Range("B28").Value = 3.272327238179451
$ws.Range("C28").Value = 1.626987699542094
$ws.Range("D28").Value = 0.1496068669990043
$ws.Range("E28").Value = 0.5333859586016987
$ws.Range("G28").Value = 5.582307763322248

$ws.Range("B29").Value = 0.2881169905109251
$ws.Range("C29").Value = 1.626987699542094
$ws.Range("D29").Value = 0.7210945179870265
$ws.Range("E29").Value = 0.5333859586016987
$ws.Range("G29").Value = 3.169585166641744

$ws.Range("B30").Value = 0.6545652718822623
$ws.Range("C30").Value = 0.04103571897497393
$ws.Range("D30").Value = 0.7210945179870265
$ws.Range("E30").Value = 0.5333859586016987
$ws.Range("G30").Value = 1.950081467445961

$ws.Range("B31").Value = 3.272327238179451
$ws.Range("C31").Value = 1.626987699542094
$ws.Range("D31").Value = 0.7210945179870265
$ws.Range("E31").Value = 0.5333859586016987
$ws.Range("G31").Value = 6.15379541431027

$ws.Range("B32").Value = 3.272327238179451
$ws.Range("C32").Value = 1.626987699542094
$ws.Range("D32").Value = 3.223369029078222
$ws.Range("E32").Value = 0.5333859586016987
$ws.Range("G32").Value = 8.656069925401464

$ws.Range("B33").Value = 0.6545652718822623
$ws.Range("C33").Value = 0.3048912486333797
$ws.Range("D33").Value = 0.1496068669990043
$ws.Range("E33").Value = 13.86384647080068
$ws.Range("G33").Value = 14.97290985831533

$ws.Range("B34").Value = 3.272327238179451
$ws.Range("C34").Value = 1.626987699542094
$ws.Range("D34").Value = 0.7210945179870265
$ws.Range("E34").Value = 0.5333859586016987
$ws.Range("G34").Value = 6.15379541431027

$ws.Range("B35").Value = 0.6545652718822623
$ws.Range("C35").Value = 1.626987699542094
$ws.Range("D35").Value = 0.1496068669990043
$ws.Range("E35").Value = 0.5333859586016987
$ws.Range("G35").Value = 2.964545797025059

$ws.Range("B36").Value = 3.272327238179451
$ws.Range("C36").Value = 1.626987699542094
$ws.Range("D36").Value = 18.71679738969934
$ws.Range("E36").Value = 0.5333859586016987
$ws.Range("G36").Value = 24.14949828602258

$ws.Range("B37").Value = 3.272327238179451
$ws.Range("C37").Value = 1.626987699542094
$ws.Range("D37").Value = 0.1496068669990043
$ws.Range("E37").Value = 0.5333859586016987
$ws.Range("G37").Value = 5.582307763322248

$ws.Range("B38").Value = 3.272327238179451
$ws.Range("C38").Value = 1.626987699542094
$ws.Range("D38").Value = 18.71679738969934
$ws.Range("E38").Value = 0.5333859586016987
$ws.Range("G38").Value = 24.14949828602258

$ws.Range("B39").Value = 3.272327238179451
$ws.Range("C39").Value = 1.626987699542094
$ws.Range("D39").Value = 3.223369029078222
$ws.Range("E39").Value = 0.5333859586016987
$ws.Range("G39").Value = 8.656069925401464

$ws.Range("B40").Value = 3.272327238179451
$ws.Range("C40").Value = 1.626987699542094
$ws.Range("D40").Value = 0.1496068669990043
$ws.Range("E40").Value = 0.5333859586016987
$ws.Range("G40").Value = 5.582307763322248

$ws.Range("B41").Value = 0.1169995834814548
$ws.Range("C41").Value = 0.3048912486333797
$ws.Range("D41").Value = 0.1496068669990043
$ws.Range("E41").Value = 0.5333859586016987
$ws.Range("G41").Value = 1.104883657715537

$ws.Range("B42").Value = 0.6545652718822623
$ws.Range("C42").Value = 1.626987699542094
$ws.Range("D42").Value = 0.1496068669990043
$ws.Range("E42").Value = 0.5333859586016987
$ws.Range("G42").Value = 2.964545797025059

$ws.Range("B43").Value = 3.272327238179451
$ws.Range("C43").Value = 1.626987699542094
$ws.Range("D43").Value = 18.71679738969934
$ws.Range("E43").Value = 0.5333859586016987
$ws.Range("G43").Value = 24.14949828602258

$ws.Range("B44").Value = 0.003078177322033415
$ws.Range("C44").Value = 0.3048912486333797
$ws.Range("D44").Value = 0.1496068669990043
$ws.Range("E44").Value = 0.5333859586016987
$ws.Range("G44").Value = 0.9909622515561161

$ws.Range("B45").Value = 1.445647641019636
$ws.Range("C45").Value = 1.626987699542094
$ws.Range("D45").Value = 18.71679738969934
$ws.Range("E45").Value = 0.5333859586016987
$ws.Range("G45").Value = 22.32281868886277

$ws.Range("B46").Value = 3.272327238179451
$ws.Range("C46").Value = 1.626987699542094
$ws.Range("D46").Value = 0.1496068669990043
$ws.Range("E46").Value = 0.5333859586016987
$ws.Range("G46").Value = 5.582307763322248

$ws.Range("B47").Value = 0.2881169905109251
$ws.Range("C47").Value = 0.3048912486333797
$ws.Range("D47").Value = 18.71679738969934
$ws.Range("E47").Value = 0.5333859586016987
$ws.Range("G47").Value = 19.84319158744534

$ws.Range("B48").Value = 3.272327238179451
$ws.Range("C48").Value = 1.626987699542094
$ws.Range("D48").Value = 0.7210945179870265
$ws.Range("E48").Value = 0.5333859586016987
$ws.Range("G48").Value = 6.15379541431027

